$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the "agustus" report row (row 9): Umum = 59, BPJS = 14
$ws.Range("B9").Value = 59
$ws.Range("C9").Value = 14

# Extend the shared SUM formula from D2 down through the new D9 row
$ws.Range("D2:D9").Formula = "=SUM(B2:C2)"

# Update the current selection shown in the workbook
$ws.Range("E10").Select()
